# Saldo_guide.xlsx — "Add files via upload" re-export
#
# The source system (IClientBalance) re-ran its export a day later
# (2024-06-24 09:40:36 -> 2024-06-25 10:15:35), so every row's reference
# date moves forward one day and the sheet is renamed to match the new
# export's timestamp. One client's balance (row 138) was also corrected
# in the refreshed extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet to the new export's timestamp
$ws.Name = "IClientBalance-20240625-101535-"

# 2. Every data row's "Dt. Referencia" (column G) advances by one day:
#    2024-06-24 (serial 45467) -> 2024-06-25 (serial 45468)
$ws.Range("G2:G277").Value = 45468

# 3. Row 138's Saldo Previsto / Vl. Total are corrected in the new extract
$ws.Range("D138").Value = 12226.59
$ws.Range("H138").Value = 12226.59

# 4. Drop the stray cell selection left over from editing (M18) so the
#    saved view reopens at the top of the sheet
$ws.Range("A1").Select()
